$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.178.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.67%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.877.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.494.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.078.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0935"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  +8.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0784"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "118.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.980.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "57.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.85%  "
